$d = $word.ActiveDocument

$replacements = @(
    @("989×8=", "639×6="),
    @("997×7=", "301×8="),
    @("528×4=", "822×4="),
    @("681×5=", "923×5="),
    @("240×6=", "924×6="),
    @("543×5=", "131×7="),
    @("947×2=", "928×7="),
    @("918×5=", "508×9="),
    @("408×9=", "334×3="),
    @("357×8=", "708×5="),
    @("929×7=", "141×9="),
    @("361×2=", "434×8="),
    @("356×2=", "780×9="),
    @("612×9=", "513×9="),
    @("319×2=", "581×7="),
    @("285×5=", "593×2="),
    @("635×5=", "311×3="),
    @("973×7=", "735×4="),
    @("413×2=", "584×5="),
    @("833×7=", "270×5="),
    @("583×3=", "828×9="),
    @("232×6=", "118×7="),
    @("121×8=", "468×6="),
    @("830×4=", "737×9="),
    @("861×5=", "706×9=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Write-Output "Done applying $($replacements.Count) replacements"
